# zxxt_02_excel_b: reshape of the demo workbook to match a resave done by
# a different user/machine (locale switched to metric page-setup units,
# sheet renamed with embedded quotes, print titles added, stale defined
# names pruned, active-sheet/selection moved to the "PrintArea" sheet).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Rename the "PrintArea (8 rows)" sheet to "PrintArea ('8' rows)" and
#    keep the Print_Area defined name (and its localSheetId) pointing at
#    the renamed sheet; add a new Print_Titles defined name for rows 2:3.
# ---------------------------------------------------------------------
$wsPrint = $wb.Worksheets.Item("PrintArea (8 rows)")
$wsPrint.Name = "PrintArea ('8' rows)"

$wb.Names.Item("PrintArea ('8' rows)!Print_Area").RefersTo = "='PrintArea (''8'' rows)'!`$A`$1:`$E`$8"
$wsPrint.PageSetup.PrintTitleRows = "`$2:`$3"

# ---------------------------------------------------------------------
# 2. Drop the stale #REF! defined names left over from earlier edits.
# ---------------------------------------------------------------------
$wb.Names.Item("E_").Delete()
$wb.Names.Item("F_").Delete()
$wb.Names.Item("RANGE_SUM1").Delete()

# ---------------------------------------------------------------------
# 3. Move the active sheet / tab selection from "ExcelTable" to the
#    renamed "PrintArea ('8' rows)" sheet, and move its selection to E1.
# ---------------------------------------------------------------------
$wsPrint.Activate()
$wsPrint.Range("E1").Select()

# ---------------------------------------------------------------------
# 4. Page setup: metric ("Normal"-under-cm) margins on both the
#    "ExcelTable" and "PrintArea ('8' rows)" sheets, plus the
#    ExcelTable sheet switching to landscape @ 80% scale.
# ---------------------------------------------------------------------
$leftRight = 0.70866141732283472 * 72
$topBottom = 0.74803149606299213 * 72
$headerFooter = 0.31496062992125984 * 72

$wsTable = $wb.Worksheets.Item("ExcelTable")
$psTable = $wsTable.PageSetup
$psTable.LeftMargin = $leftRight
$psTable.RightMargin = $leftRight
$psTable.TopMargin = $topBottom
$psTable.BottomMargin = $topBottom
$psTable.HeaderMargin = $headerFooter
$psTable.FooterMargin = $headerFooter
$psTable.Orientation = 2
$psTable.Zoom = 80

$psPrint = $wsPrint.PageSetup
$psPrint.LeftMargin = $leftRight
$psPrint.RightMargin = $leftRight
$psPrint.TopMargin = $topBottom
$psPrint.BottomMargin = $topBottom
$psPrint.HeaderMargin = $headerFooter
$psPrint.FooterMargin = $headerFooter
